$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 11, shifting all existing data (rows 11-64) down to rows 13-66.
$ws.Rows.Item(11).Resize(2).Insert()

# Populate the new row 11 with the latest "Especial" quality observation.
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C11").Value = "Metropolitana"
$ws.Range("D11").Value = 45071
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100107
$ws.Range("H11").Value = "Otros"
$ws.Range("I11").Value = 100107001
$ws.Range("J11").Value = "Caqui"
$ws.Range("K11").Value = "Mankaki"
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 330
$ws.Range("N11").Value = 15000
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 15000
$ws.Range("Q11").Value = "`$/caja 16 kilos granel"
$ws.Range("R11").Value = "Región de O'Higgins"
$ws.Range("S11").Value = 938
$ws.Range("T11").Value = 16

# Populate the new row 12 with the latest "Primera" quality observation.
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C12").Value = "Metropolitana"
$ws.Range("D12").Value = 45071
$ws.Range("E12").Value = 13
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100107
$ws.Range("H12").Value = "Otros"
$ws.Range("I12").Value = 100107001
$ws.Range("J12").Value = "Caqui"
$ws.Range("K12").Value = "Mankaki"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 280
$ws.Range("N12").Value = 13000
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 13000
$ws.Range("Q12").Value = "`$/caja 16 kilos granel"
$ws.Range("R12").Value = "Región de O'Higgins"
$ws.Range("S12").Value = 812
$ws.Range("T12").Value = 16
